$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "展览" (Worksheets index 1)
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# F2: 1019 -> 1020
$ws1.Range("F2").Value = 1020

# F3: 2109 -> 2123
$ws1.Range("F3").Value = 2123

# Insert a new row above the current row 4 (old row 4 "布谷鸟动漫展4th"
# shifts down to row 5, keeping its values/format).
$ws1.Rows.Item(4).Insert()

# Copy the number-column format (bold/border/center) from row 3's A cell
# onto the newly inserted A4 so it matches the rest of the column.
$ws1.Range("A3").Copy()
$ws1.Range("A4").PasteSpecial(-4122)

# Populate the freshly inserted row 4 with the new event.
$ws1.Range("A4").Value = 3
# Leading apostrophe keeps this looking like a date string from being
# auto-converted into a date serial number; re-applying the "Normal"
# style afterwards drops the quote-prefix formatting flag again so the
# cell ends up styled exactly like its neighbours.
$ws1.Range("B4").Value = "'2024-05-19"
$ws1.Range("B4").Style = "Normal"
$ws1.Range("C4").Value = "南宁·原x穹x崩only"
$ws1.Range("D4").Value = "明秀东路157号 利泰国际大酒店"
$ws1.Range("E4").Value = "2024.05.19 10:00-05.19 17:00"
$ws1.Range("F4").Value = 0
$ws1.Range("G4").Value = 35
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=83070"
$ws1.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

# Row 5 now holds the event that used to live in row 4; fix its index
# number and the updated "want to go" count.
$ws1.Range("A5").Value = 4
$ws1.Range("F5").Value = 464

# ----------------------------------------------------------------------
# Sheet "全部类型" (Worksheets index 4)
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# F4: 1019 -> 1020
$ws4.Range("F4").Value = 1020

# F5: 2109 -> 2123
$ws4.Range("F5").Value = 2123

# Insert a new row above the current row 6 (old row 6 "布谷鸟动漫展4th"
# shifts down to row 7, keeping its values/format).
$ws4.Rows.Item(6).Insert()

# Copy the number-column format from row 5's A cell onto the newly
# inserted A6 so it matches the rest of the column.
$ws4.Range("A5").Copy()
$ws4.Range("A6").PasteSpecial(-4122)

# Populate the freshly inserted row 6 with the new event.
$ws4.Range("A6").Value = 5
$ws4.Range("B6").Value = "'2024-05-19"
$ws4.Range("B6").Style = "Normal"
$ws4.Range("C6").Value = "南宁·原x穹x崩only"
$ws4.Range("D6").Value = "明秀东路157号 利泰国际大酒店"
$ws4.Range("E6").Value = "2024.05.19 10:00-05.19 17:00"
$ws4.Range("F6").Value = 0
$ws4.Range("G6").Value = 35
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=83070"
$ws4.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202403/I8tScigE1710918412731.jpeg"

# Row 7 now holds the event that used to live in row 6; fix its index
# number and the updated "want to go" count.
$ws4.Range("A7").Value = 6
$ws4.Range("F7").Value = 464
